# control de coma por punto decimal
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the xn values to use dot-decimal / semicolon-separated vector formatting
$ws.Range("B2").Value = "[1.6;1.75;1]"
$ws.Range("B3").Value = "[1.6;1.75;1]"

# The third iteration row is no longer needed; its error value (0) now belongs to row 3.
# Force text formatting so the numeric-looking value stays a text string (matches the
# rest of the sheet, which stores every value - including numbers - as text).
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = "0"

# Remove the now-obsolete row 4 entirely (shrinks the used range to A1:C3)
$ws.Rows("4:4").Delete()
